$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 21817
$ws.Range("G3").Value = 95
$ws.Range("F5").Value = 8153
$ws.Range("F6").Value = 568
$ws.Range("F7").Value = 333
$ws.Range("F8").Value = 75
$ws.Range("F9").Value = 209
$ws.Range("F13").Value = 252
$ws.Range("F14").Value = 1037
$ws.Range("F15").Value = 1377
$ws.Range("F18").Value = 725
$ws.Range("F20").Value = 109
$ws.Range("F22").Value = 380
$ws.Range("F23").Value = 1246
$ws.Range("F26").Value = 249
$ws.Range("F27").Value = 5213
$ws.Range("F30").Value = 177
$ws.Range("F31").Value = 5306
$ws.Range("F32").Value = 41
$ws.Range("F36").Value = 13610
$ws.Range("F38").Value = 170
$ws.Range("F39").Value = 68
$ws.Range("F41").Value = 359
$ws.Range("F42").Value = 500
$ws.Range("F43").Value = 4116
$ws.Range("F44").Value = 59
$ws.Range("F45").Value = 339

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 39

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 21817
$ws.Range("G3").Value = 95
$ws.Range("F5").Value = 8153
$ws.Range("F6").Value = 568
$ws.Range("F7").Value = 333
$ws.Range("F8").Value = 75
$ws.Range("F9").Value = 209
$ws.Range("F13").Value = 252
$ws.Range("F14").Value = 1037
$ws.Range("F15").Value = 1377
$ws.Range("F18").Value = 725
$ws.Range("F20").Value = 109
$ws.Range("F22").Value = 380
$ws.Range("F23").Value = 1246
$ws.Range("F26").Value = 249
$ws.Range("F28").Value = 5213
$ws.Range("F32").Value = 177
$ws.Range("F33").Value = 39
$ws.Range("F34").Value = 5306
$ws.Range("F35").Value = 41
$ws.Range("F39").Value = 13610
$ws.Range("F41").Value = 170
$ws.Range("F42").Value = 68
$ws.Range("F44").Value = 359
$ws.Range("F45").Value = 500
$ws.Range("F46").Value = 4116
$ws.Range("F47").Value = 59
$ws.Range("F48").Value = 339
